$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All rows 2-51 have column G (Hora) change from "4" to "5".
# Apply Text format to the whole contiguous G2:G51 range once, then write values.
$ws.Range("G2:G51").NumberFormat = "@"
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 7).Value = "5"
}

# Column B/C/E are plain (non-numeric-looking) text -> no NumberFormat needed.
# Column D values look numeric, so force Text format before writing so Excel
# keeps them as text (matching the source workbook which stores them as
# inline strings, not numbers).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "248.82"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.74"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.367"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05616"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.410"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.375"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8154"

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9613"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1412"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07572"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03178"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03051"

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09312"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.554"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001604"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04714"

# Row 18
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005774"
$ws.Range("E18").Value = "17OneONEWorstin24h"

# Row 19
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006423"
$ws.Range("E19").Value = "18TigerCashTCH"

# Row 20
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005087"
$ws.Range("E20").Value = "19HotbitTokenHTB"

# Row 21
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001033"
$ws.Range("E21").Value = "20BitKanKAN"

# Row 22
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001499"
$ws.Range("E22").Value = "21NitroExNTX"

# Row 23
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.746"
$ws.Range("E23").Value = "22LEOLEO"

# Row 24
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.142"
$ws.Range("E24").Value = "23BTSETokenBTSE"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3255"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1250"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006964"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1064"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002997"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008619"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0005495"
$ws.Range("E47").Value = "46ACDXExchangeACXT"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.6794"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1610"
$ws.Range("E49").Value = "48BOLOBOLOBestin24h"
